$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update existing value cells in row 3
$ws.Range("H3").Value = -5.7
$ws.Range("I3").Value = 4.9000000000000004
$ws.Range("N3").Value = -5.7
$ws.Range("O3").Value = 6.2
$ws.Range("P3").Value = 0.9
$ws.Range("Q3").Value = 1.8

# Clear cells that were removed (J3:M3)
$ws.Range("J3:M3").ClearContents()

# Update the selection to reflect the new active cell/selection
$ws.Range("M3").Select()
